# Natmi following Dr Hou advice
# Rewrite the LR-pair result rows (2-10) with updated sender/receiver clusters
# and refreshed expression-specificity statistics, per the revised NATMI run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("ECs", "Fbln1", "Itgb1", "ECs", 1, 0.3333333333333333, 0.08550066666666667, 0.256502, 0.001136245970483483, 0.001136245970483483, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 9.619970024928, 86.57973022435199, 0.0003721561172607291, 0.0003721561172607291),
    @("ECs", "Fbln1", "Itgb1", "FAPs", 1, 0.3333333333333333, 0.08550066666666667, 0.256502, 0.001136245970483483, 0.001136245970483483, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 9.089957719310666, 81.809619473796, 0.0003516521737715257, 0.0003516521737715258),
    @("ECs", "Fbln1", "Itgb1", "sCs", 1, 0.3333333333333333, 0.08550066666666667, 0.256502, 0.001136245970483483, 0.001136245970483483, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 10.66121965877022, 95.95097692893201, 0.0004124376794512285, 0.0004124376794512286),
    @("FAPs", "Fbln1", "Itgb1", "ECs", 3, 1, 72.45391466666668, 217.361744, 0.9628634699037922, 0.9628634699037921, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 8152.035702825217, 73368.32132542695, 0.3153679218410016, 0.3153679218410015),
    @("FAPs", "Fbln1", "Itgb1", "FAPs", 3, 1, 72.45391466666668, 217.361744, 0.9628634699037922, 0.9628634699037921, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 7702.899247396235, 69326.09322656611, 0.2979927243154826, 0.2979927243154826),
    @("FAPs", "Fbln1", "Itgb1", "sCs", 3, 1, 72.45391466666668, 217.361744, 0.9628634699037922, 0.9628634699037921, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 9034.398555166747, 81309.58699650072, 0.349502823747308, 0.3495028237473081),
    @("sCs", "Fbln1", "Itgb1", "ECs", 3, 1, 2.708963, 8.126889, 0.03600028412572435, 0.03600028412572434, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 304.794615932496, 2743.151543392464, 0.01179121977859405, 0.01179121977859404),
    @("sCs", "Fbln1", "Itgb1", "FAPs", 3, 1, 2.708963, 8.126889, 0.03600028412572435, 0.03600028412572434, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 288.001954758758, 2592.017592828822, 0.01114158245491225, 0.01114158245491225),
    @("sCs", "Fbln1", "Itgb1", "sCs", 3, 1, 2.708963, 8.126889, 0.03600028412572435, 0.03600028412572434, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 337.7850807067528, 3040.065726360774, 0.01306748189221805, 0.01306748189221805)
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = 2 + $i
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
